# Update "想去人数" (F2/F3) figures on the "展览" and "全部类型" sheets
# F2: 129 -> 130
# F3: 51  -> 53

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 130
    $ws.Range("F3").Value = 53
}
